# Apply the edit described by the diff:
#  - Row 2 data is replaced with new strategy results (ticker SEMI.BA, updated metrics)
#  - Rows 3-5 (MORI.BA, EDN.BA, YPFD.BA) are removed entirely
#  - Sheet dimension shrinks from A1:W5 to A1:W2 (handled automatically once rows are deleted)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new values
$ws.Range("A2").Value = "SEMI.BA"
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 10000
$ws.Range("E2").Value = 65754.37359242247
$ws.Range("G2").Value = 45170
$ws.Range("H2").Value = 42
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 0.333
$ws.Range("L2").Value = 0.667
$ws.Range("M2").Value = -0.08699999999999999
$ws.Range("N2").Value = 5709.358
$ws.Range("O2").Value = 0.271
$ws.Range("P2").Value = -1456.424
$ws.Range("Q2").Value = -0.06
$ws.Range("R2").Value = -6238.426
$ws.Range("S2").Value = 45693.333
$ws.Range("T2").Value = 4.569333321689681
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 7
$ws.Range("W2").Value = 34.26829268292683

# Remove rows 3 through 5 (MORI.BA, EDN.BA, YPFD.BA strategies)
$ws.Rows("3:5").Delete()
